$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new sales-order row (row 7) that was entered for Ashley Baker /
# Korn / "Whatever".
#
# A7 ("01/12/2024") and D7 ("9157994875") look numeric/date-like to Excel's
# auto-detection, so they are forced to Text via NumberFormat "@" before the
# value is written - this is exactly the PHONE-field-shown-as-a-float bug
# called out in the commit message, now fixed for the new row too. The
# style is reset back to "Normal" right after so the cell doesn't end up
# carrying an explicit style index - same as the other (unstyled) data rows,
# which simply inherit their column's default style.

function Set-TextCell($row, $col, $value) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.Value = $value
    $cell.Style = "Normal"
}

function Set-ForcedTextCell($row, $col, $value) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

function Set-NumberCell($row, $col, $value) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-ForcedTextCell 7 1 "01/12/2024"
Set-TextCell       7 2 "SO240112001"
Set-TextCell       7 3 "Ashley Baker"
Set-ForcedTextCell 7 4 "9157994875"
Set-TextCell       7 5 "Korn"
Set-TextCell       7 6 "Whatever"
Set-NumberCell     7 7 50
Set-NumberCell     7 8 95
Set-TextCell       7 9 "AMS"
Set-TextCell       7 10 "DVD"
Set-TextCell       7 11 "Ashley"
